{"js": "// The edit: \" vous vous en trouver\u00e9s fort bien<lb/>.\" (spread across three\n// runs \u2014 text, a gray \"<lb/>\" line-break marker, and a final \".\") collapses\n// into a single run reading \" vous vous en trouver\u00e9s fort bien.\" \u2014 i.e. the\n// \"<lb/>\" marker is removed and the trailing period is pulled up onto the\n// same run/formatting as the preceding sentence.\nconst body = context.document.body;\n\n// Locate the run of text spanning \"...fort bien<lb/>.\" \u2014 unique in the doc.\nconst needle = \"vous vous en trouver\u00e9s fort bien<lb/>.\";\nconst results = body.search(needle, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\n    `expected exactly 1 match for \"${needle}\", found ${results.items.length}`\n  );\n}\n\nconst target = results.items[0];\n// Replace the whole matched span with the merged text; Word keeps the\n// formatting of the first run in the replaced range (the plain, non-\"<lb/>\"\n// run), which is exactly the run the final text should carry.\ntarget.insertText(\"vous vous en trouver\u00e9s fort bien.\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The edit: \" vous vous en trouver\u00e9s fort bien<lb/>.\" (spread across three\n# runs \u2014 text, a gray \"<lb/>\" line-break marker, and a final \".\") collapses\n# into a single run reading \" vous vous en trouver\u00e9s fort bien.\" \u2014 i.e. the\n# \"<lb/>\" marker is removed and the trailing period is pulled up right after\n# \"bien\" with no marker between them.\n\n$d = $word.ActiveDocument\n\n$needleFind    = \"vous vous en trouver\u00e9s fort bien<lb/>.\"\n$needleReplace = \"vous vous en trouver\u00e9s fort bien.\"\n\n# Sanity-check: exactly one occurrence in the document before editing.\n$before = $d.Content.Text\n$hitCount = ([regex]::Matches($before, [regex]::Escape($needleFind))).Count\nif ($hitCount -ne 1) {\n    throw \"expected exactly 1 match for '$needleFind', found $hitCount\"\n}\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $needleFind\n$find.Replacement.Text = $needleReplace\n\n# wdFindContinue (1) keeps the search scoped forward without prompting;\n# wdReplaceAll (2) rewrites every match (there is exactly one, per the check\n# above).\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\nif (-not $found) {\n    throw \"Find.Execute did not find '$needleFind'\"\n}\n\n# Verify the replacement actually landed and the marker is gone.\n$after = $d.Content.Text\nif ($after -notlike \"*$needleReplace*\") {\n    throw \"replacement text not found after Execute\"\n}\nif ($after -like \"*<lb/>.*\") {\n    throw \"<lb/>. marker still present after Execute\"\n}\n"}
